# Generate Report for Handoff
#
# - Overview sheet: status text for the handed-off file flips from
#   "Handoff transform failed" to "Ready for handoff".
# - zh-cn / de-de sheets: the handoff that used to be blocked now has a
#   real handoff file + timestamp recorded, and the "Handoff Reason"
#   flips from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2          # xlUnderlineStyleSingle
$hyperlinkColor = 15570276       # RGB(0x64, 0x95, 0xED) == cornflower blue, BGR-packed

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhFile = "d473859f-d9d0-42d3-96c3-ffddf165fab4.b902c4b270c87a80eef8d92639f106c98599481e.zh-cn.xlf"
$zhUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c897c10fd19a949f0b99fdad03d405078b51ee63/" + $zhFile

$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "2016-02-17 06:32:48"
$wsZh.Range("H2").Value = "Include"

# Rebuild the hyperlinks in final left-to-right / top-to-bottom order so the
# new handoff-file link lands between the existing two (matching how Excel
# renumbers r:id on save).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c897c10fd19a949f0b99fdad03d405078b51ee63/e2e/d473859f-d9d0-42d3-96c3-ffddf165fab4.md", "", "", "d473859f-d9d0-42d3-96c3-ffddf165fab4.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhUrl, "", "", $zhFile)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c897c10fd19a949f0b99fdad03d405078b51ee63/.localization-config", "", "", ".localization-config")

$wsZh.Range("A2").Font.Underline = $hyperlinkUnderline
$wsZh.Range("A2").Font.Color = $hyperlinkColor
$wsZh.Range("C2").Font.Underline = $hyperlinkUnderline
$wsZh.Range("C2").Font.Color = $hyperlinkColor
$wsZh.Range("A3").Font.Underline = $hyperlinkUnderline
$wsZh.Range("A3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deFile = "d473859f-d9d0-42d3-96c3-ffddf165fab4.b902c4b270c87a80eef8d92639f106c98599481e.de-de.xlf"
$deUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c897c10fd19a949f0b99fdad03d405078b51ee63/" + $deFile

$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "2016-02-17 06:32:59"
$wsDe.Range("H2").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c897c10fd19a949f0b99fdad03d405078b51ee63/e2e/d473859f-d9d0-42d3-96c3-ffddf165fab4.md", "", "", "d473859f-d9d0-42d3-96c3-ffddf165fab4.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deUrl, "", "", $deFile)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c897c10fd19a949f0b99fdad03d405078b51ee63/.localization-config", "", "", ".localization-config")

$wsDe.Range("A2").Font.Underline = $hyperlinkUnderline
$wsDe.Range("A2").Font.Color = $hyperlinkColor
$wsDe.Range("C2").Font.Underline = $hyperlinkUnderline
$wsDe.Range("C2").Font.Color = $hyperlinkColor
$wsDe.Range("A3").Font.Underline = $hyperlinkUnderline
$wsDe.Range("A3").Font.Color = $hyperlinkColor

Write-Output "Handoff report generated"
